# Updates the cryptocurrency Price (D) and Volume/1h-change (E) columns
# for rows 2-51 on Sheet1, matching the latest GitHub Actions data refresh.
#
# Note: several "Price" values are plain decimal numbers (e.g. "1.001",
# "74.00"). Excel's Range.Value setter auto-converts such numeric-looking
# strings to real numbers (losing exact text / trailing zeros), exactly as
# it would if a user typed them into a General-formatted cell. To preserve
# them as literal text - as they are stored in the workbook - we prefix
# those values with a leading apostrophe, which is the standard Excel
# "treat as text" input convention (sets the quote-prefix flag) instead of
# reformatting the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.911.04"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.895.91"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'0.7833"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'244.15"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.3144"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").Value = "'25.73"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "'0.07274"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").Value = "'0.08119"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "'0.7741"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "'5.479"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "1.893.41"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'94.39"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").Value = "'6.218"
$ws.Range("E16").Value = "  +5.10%  "
$ws.Range("D17").Value = "29.896.55"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "'13.95"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'246.16"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'0.000007833"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'8.132"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "2.122.58"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'0.1606"
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("D26").Value = "'9.467"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "'164.46"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").Value = "'1.434"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'4.478"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").Value = "'0.05582"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "'4.084"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "'0.7540"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").Value = "'0.9968"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'2.679"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").Value = "'0.01936"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "'2.789"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "1.142.51"
$ws.Range("E41").Value = "  +11.86%  "
$ws.Range("D42").Value = "'0.4461"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'74.00"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").Value = "'5.966"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").Value = "'0.8541"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'1.891"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "'3.148"
$ws.Range("E48").Value = "  +8.62%  "
$ws.Range("D49").Value = "'102.06"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "'9.797"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'7.534"
$ws.Range("E51").Value = "  +1.74%  "
